{"js": "// The document's visible text is unchanged (the bulk of the underlying\n// diff is just Word's own proofing engine splitting runs and tagging\n// them with <w:proofErr> spell/grammar markers - invisible to the\n// reader and not something an editing script should try to recreate).\n// The one real, content-visible change is a new empty paragraph added\n// at the very end of the body, right after the \"npm start\" line,\n// inheriting that line's run formatting (sz/szCs = 40).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Append a new, empty paragraph after the last one in the document.\nconst newParagraph = lastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\n\n// Keep the same font size as the rest of the document (OOXML w:sz/w:szCs\n// are in half-points, so sz=\"40\" == 20pt) so the inserted paragraph mark\n// carries the same run properties seen throughout the file.\nnewParagraph.font.size = 20;\n\nawait context.sync();\n", "ps1": "# The visible text of the document is unchanged by this revision (the\n# bulk of the underlying diff is just Word's own proofing engine\n# splitting runs apart and wrapping pieces in <w:proofErr> spell/grammar\n# markers - invisible to the reader, and not something a script should\n# try to recreate). The one real, content-visible change is a new empty\n# paragraph added at the very end of the document, right after the\n# \"npm start\" line, inheriting that line's run formatting (sz/szCs = 40,\n# i.e. 20pt).\n\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n$newLastParagraph = $d.Paragraphs.Last\n$newLastParagraph.Range.Font.Size = 20\n"}
